# Moved some BoM components to Assembly BoM
#
# The "Need contacts, wire, mating connectors, terminating connectors" note
# (row 16) is no longer needed on this BoM sheet, so remove the entire row.
# Deleting the row shifts the Subtotal row (and the two blank spacer rows
# below it) up by one, which matches the target layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select row 16 first (as a user would before deleting it) so the resulting
# selection/active cell state matches what Excel records after the edit.
$ws.Rows(16).Select()

# Remove the whole row - this shifts everything below it up by one row,
# automatically fixing the dimension, formulas, and hyperlink references.
$ws.Rows(16).Delete()
